$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 3.182878228561681
$ws.Range("C2").Value = 1.65323645889881
$ws.Range("D2").Value = 0.1529057820181812
$ws.Range("E2").Value = 6.48142807727062
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 11.47044854674929

# Row 3
$ws.Range("B3").Value = 1.505614041169197
$ws.Range("C3").Value = 0.3375848360084654
$ws.Range("D3").Value = 0.1529057820181812
$ws.Range("E3").Value = 0.4998867070740569
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 2.495991366269901
